# Generate Report for Handoff
# Replace the file that was handed off (old GUID) with the newly generated
# handoff file (new GUID), refresh the handoff/handback timestamps, and
# clear out the stale "target"/"handback" info for zh-cn and de-de since a
# fresh handoff was just generated (no handback has happened for it yet).

$wb = $excel.ActiveWorkbook

$newGuid = "75c54032-39f2-4fac-b439-864bc2ddc7d1"

# ============================================================
# Sheet: Overview
# ============================================================
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
foreach ($h in $wsOverview.Hyperlinks) {
    if ($h.Range.Column -eq 2) {
        $h.TextToDisplay = "e2e\$newGuid.md"
    }
}
$wsOverview.Range("G2").Value = "2016-09-04 23:04:54"

# ============================================================
# Sheet: zh-cn
# ============================================================
$wsZh = $wb.Worksheets.Item("zh-cn")

# Drop the hyperlink that lived on the (now stale) "Latest Target File"
# cell, but keep the one on "Source File Name" (column A).
foreach ($h in $wsZh.Hyperlinks) {
    if ($h.Range.Column -eq 9) {
        $h.Delete()
    }
}
foreach ($h in $wsZh.Hyperlinks) {
    if ($h.Range.Column -eq 1) {
        $h.TextToDisplay = "$newGuid.md"
    }
}

$wsZh.Range("A2").Value = "$newGuid.md"
$wsZh.Range("H2").Value = "2016-09-04 23:04:50"
$wsZh.Range("I2").Value = ""
$wsZh.Range("I2").Style = "Normal"
$wsZh.Range("J2").Value = ""
$wsZh.Range("K2").Value = "0001-01-01 00:00:00"

$wsZh.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsZh.Columns.Item(10).ColumnWidth = 21.7054770333426

# ============================================================
# Sheet: de-de
# ============================================================
$wsDe = $wb.Worksheets.Item("de-de")

foreach ($h in $wsDe.Hyperlinks) {
    if ($h.Range.Column -eq 9) {
        $h.Delete()
    }
}
foreach ($h in $wsDe.Hyperlinks) {
    if ($h.Range.Column -eq 1) {
        $h.TextToDisplay = "$newGuid.md"
    }
}

$wsDe.Range("A2").Value = "$newGuid.md"
$wsDe.Range("H2").Value = "2016-09-04 23:04:54"
$wsDe.Range("I2").Value = ""
$wsDe.Range("I2").Style = "Normal"
$wsDe.Range("J2").Value = ""
$wsDe.Range("K2").Value = "0001-01-01 00:00:00"

$wsDe.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsDe.Columns.Item(10).ColumnWidth = 21.7054770333426
